$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '26.082.19'
Set-TextValue 'E2' '  -1.71%  '
Set-TextValue 'D3' '1.665.72'
Set-TextValue 'E3' '  -1.17%  '
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '216.31'
Set-TextValue 'E5' '  -0.07%  '
Set-TextValue 'D6' '0.5106'
Set-TextValue 'E6' '  +2.24%  '
Set-TextValue 'E7' '  -0.01%  '
Set-TextValue 'E8' '  +0.85%  '
Set-TextValue 'D9' '0.06423'
Set-TextValue 'E9' '  +4.47%  '
Set-TextValue 'D10' '21.55'
Set-TextValue 'E10' '  -0.66%  '
Set-TextValue 'D11' '0.07418'
Set-TextValue 'E11' '  +1.75%  '
Set-TextValue 'D12' '1.665.58'
Set-TextValue 'E12' '  -1.13%  '
Set-TextValue 'D13' '4.509'
Set-TextValue 'E13' '  +1.79%  '
Set-TextValue 'D14' '0.5800'
Set-TextValue 'E14' '  +1.20%  '
Set-TextValue 'D15' '0.000008555'
Set-TextValue 'E15' '  +3.49%  '
Set-TextValue 'D16' '64.23'
Set-TextValue 'E16' '  -0.87%  '
Set-TextValue 'D17' '26.129.95'
Set-TextValue 'E17' '  -1.70%  '
Set-TextValue 'D18' '4.917'
Set-TextValue 'E18' '  -2.03%  '
Set-TextValue 'E19' '  -0.07%  '
Set-TextValue 'E20' '  +0.69%  '
Set-TextValue 'D21' '189.29'
Set-TextValue 'E21' '  +3.28%  '
Set-TextValue 'D22' '6.208'
Set-TextValue 'E22' '  +0.39%  '
Set-TextValue 'D23' '1.006'
Set-TextValue 'E23' '  -0.01%  '
Set-TextValue 'D24' '145.30'
Set-TextValue 'E24' '  +0.43%  '
Set-TextValue 'D25' '7.615'
Set-TextValue 'E25' '  +0.46%  '
Set-TextValue 'D26' '0.1201'
Set-TextValue 'E26' '  +5.93%  '
Set-TextValue 'E27' '  +1.67%  '
Set-TextValue 'D28' '0.06402'
Set-TextValue 'E28' '  +14.30%  '
Set-TextValue 'D29' '1.297'
Set-TextValue 'E29' '  -1.49%  '
Set-TextValue 'D30' '1.315'
Set-TextValue 'E30' '  -0.56%  '
Set-TextValue 'D31' '3.523'
Set-TextValue 'E31' '  +1.34%  '
Set-TextValue 'D32' '3.504'
Set-TextValue 'E32' '  +0.94%  '
Set-TextValue 'E33' '  -0.29%  '
Set-TextValue 'D34' '1.014'
Set-TextValue 'E34' '  +0.79%  '
Set-TextValue 'D35' '0.6087'
Set-TextValue 'E35' '  +3.48%  '
Set-TextValue 'D36' '2.360'
Set-TextValue 'E36' '  -0.51%  '
Set-TextValue 'D37' '2.649'
Set-TextValue 'E37' '  +0.47%  '
Set-TextValue 'D38' '6.155'
Set-TextValue 'E38' '  +4.11%  '
Set-TextValue 'D39' '0.01608'
Set-TextValue 'E39' '  +0.94%  '
Set-TextValue 'D40' '1.075.75'
Set-TextValue 'E40' '  -0.07%  '
Set-TextValue 'D41' '0.8599'
Set-TextValue 'E41' '  +0.64%  '
Set-TextValue 'E42' '  +0.66%  '
Set-TextValue 'E43' '  +2.51%  '
Set-TextValue 'D44' '1.814.21'
Set-TextValue 'E44' '  -1.53%  '
Set-TextValue 'E45' '  +8.93%  '
Set-TextValue 'E46' '  -0.34%  '
Set-TextValue 'D47' '1.007'
Set-TextValue 'E47' '  +0.27%  '
Set-TextValue 'D48' '8.093'
Set-TextValue 'E48' '  +0.24%  '
Set-TextValue 'D49' '0.05205'
Set-TextValue 'E49' '  +0.04%  '
Set-TextValue 'D50' '0.4287'
Set-TextValue 'E50' '  -0.95%  '
Set-TextValue 'E51' '  +6.37%  '
